# Tweaks to baseline data: set "MFG Eff" column B (Average Module Efficiency)
# values to reflect Willeke's module MFG eff of 98% in 2002 -> ramping from
# 95% up to 99% by 2030 and held flat afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MFG Eff")

$values = @(
    95, 95, 95, 95, 95, 95, 95, 95, 95, 95,
    95, 95, 95, 95, 95, 95, 95, 95, 95, 95.2,
    95.5, 95.7, 95.9, 96.2, 96.4, 96.6, 96.9, 97.1, 97.4, 97.6,
    97.8, 98.1, 98.3, 98.5, 98.8, 99, 99, 99, 99, 99,
    99, 99, 99, 99, 99, 99, 99, 99, 99, 99,
    99, 99, 99, 99, 99, 99
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("J31").Select()
